$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values (columns B, C, D, E, G, I, M, O) for rows 2-25
$colB = @(0.8204252994459011,0.7209362687368639,0.659743177468215,0.6347808497642973,0.6306343585291643,0.6594066291340823,0.7861442032865966,1.03379563686542,1.215178727098646,1.297566782645333,1.32874637193521,1.322032151640656,1.30013233308415,1.286715556559386,1.209791879939587,1.162569109938147,1.135396215309584,1.126193981896336,1.167597263990501,1.306565364775111,1.397277908690342,1.348873459713275,1.165324109893731,0.9668967683517735)
$colC = @(0.227809183241817,0.1988024612056734,0.180927734413217,0.1736276390377895,0.172414503921118,0.1808293472631419,0.2178211679961635,0.2898445855980469,0.3424423937881897,0.3663014075121964,0.3753262981038574,0.3733830746868989,0.3670440923673368,0.3631599738315572,0.3408817637737229,0.3271972286641471,0.3193198591796431,0.3166516271508613,0.3286546334348088,0.3689062782827932,0.3951546671715391,0.3811508227670402,0.3279957717134891,0.2704160010221131)
$colD = @(0.07914228364920461,0.07169747977485486,0.06716268582437124,0.06532380698453721,0.06501900982534892,0.06713784933897671,0.07656773746202816,0.09535148505619873,0.1093358786027494,0.1157390900858104,0.1181698889502769,0.1176461033605847,0.1159389519786487,0.1148940608650122,0.1089182561827471,0.1052629969336323,0.1031645046907244,0.1024546622489737,0.1056516998212231,0.1164402193235929,0.1235263809191451,0.1197411207872392,0.1054759580170952,0.09023811779591995)
$colE = @(0.1058171064006501,0.1008576460051742,0.09792283209908703,0.09675429353514886,0.09656190554123256,0.09790696216793293,0.1040839955724238,0.1170874768383143,0.1272057529775381,0.1319362027396167,0.1337462170739911,0.1333555610510047,0.1320847370209464,0.1313087656870238,0.1268992041766595,0.1242270079671073,0.1227020203565417,0.1221877343088522,0.1245102244909901,0.1324574986686855,0.1377606371543791,0.134920151079271,0.1243821471510813,0.1134721220575869)
$colG = @(0.002403682274624849,0.002407052235307351,0.002409230710843558,0.002410146030363588,0.002410299686342405,0.002409242943503598,0.002404821603029295,0.002397014684295657,0.002391799659909136,0.002389539089396481,0.002388699052317359,0.002388879259232144,0.002389469658962862,0.002389833376094343,0.002391949635708685,0.002393276462245508,0.002394050142838516,0.002394313907917706,0.002393134130543738,0.002389295811038158,0.00238688042318417,0.002388161062054737,0.002393198444703954,0.00239903481776811)
$colI = @(0.4260496289937556,0.4319414183826495,0.4359679241284624,0.4377112078998096,0.438006855519717,0.4359910203498849,0.4279960193651249,0.4155794524748444,0.4084688361303392,0.4056759397182859,0.4046822819518141,0.4048934328649594,0.4055929067650226,0.4060296957799565,0.408660284297607,0.4103875139541167,0.4114225256193258,0.4117800868652282,0.410199343541283,0.4053857152939599,0.4026126876815894,0.4040584379239576,0.4102842845713788,0.418586860113848)
$colM = @(0.3569098304336933,0.3186550151146577,0.295237403581865,0.2857123006725928,0.2841317350424077,0.2951088728990285,0.3437047666807445,0.4395737624336249,0.5103780871567807,0.542673921599885,0.5549162098779021,0.552279053856779,0.5436808497466501,0.5384158399314742,0.508269229697774,0.4897975116971622,0.4791812118659351,0.4755881107647753,0.4917630140234621,0.5462060090659833,0.5818608966495162,0.5628244740981359,0.4908744000441487,0.4135752828674555)
$colO = @(2.163602134772987,2.155727056313566,2.152940309686073,2.152317686022855,2.152245221183279,2.152929838497244,2.160460330315573,2.19158285586596,2.224564832008838,2.241797582237496,2.248646164732435,2.247156800430787,2.242354533955563,2.239455130517797,2.223483689042098,2.214258337662159,2.209161829365172,2.207472177922398,2.215218674821614,2.243756293299498,2.264290570200387,2.253157933521379,2.214783860923063,2.181397396028132)

$rows = 2..25
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $ws.Range("B$r").Value = $colB[$i]
    $ws.Range("C$r").Value = $colC[$i]
    $ws.Range("D$r").Value = $colD[$i]
    $ws.Range("E$r").Value = $colE[$i]
    $ws.Range("G$r").Value = $colG[$i]
    $ws.Range("I$r").Value = $colI[$i]
    $ws.Range("M$r").Value = $colM[$i]
    $ws.Range("O$r").Value = $colO[$i]
}
